# Changes as per coach request
$wb = $excel.ActiveWorkbook

# --- "category" sheet: replace the category list ---
$wsCategory = $wb.Worksheets.Item("category")

# Clear the old content range (A1:B10) then write the new, shorter table
$wsCategory.Range("A1:B10").ClearContents()

$wsCategory.Range("A1").Value = "Category"

$wsCategory.Range("A2").Value = "Carro"
$wsCategory.Range("B2").Value = "Expense"

$wsCategory.Range("A3").Value = "Universidad"
$wsCategory.Range("B3").Value = "Expense"

$wsCategory.Range("A4").Value = "Servicios Profesionales"
$wsCategory.Range("B4").Value = "Income"

# --- "Expense and incomes" sheet: trim down to two transactions ---
$wsExpense = $wb.Worksheets.Item("Expense and incomes")

$wsExpense.Range("A1:D7").ClearContents()

$wsExpense.Range("A1").Value = "Detail"
$wsExpense.Range("B1").Value = "Category"
$wsExpense.Range("C1").Value = "Type"
$wsExpense.Range("D1").Value = "Amount"

$wsExpense.Range("A2").Value = "Matricula"
$wsExpense.Range("B2").Value = "Universidad"
$wsExpense.Range("C2").Value = "Expense"
$wsExpense.Range("D2").Value = "'4500"

$wsExpense.Range("A3").Value = "Analisis de un sistema informatico"
$wsExpense.Range("B3").Value = "Servicios Profesionales"
$wsExpense.Range("C3").Value = "Income"
$wsExpense.Range("D3").Value = "'45000"

# The leading apostrophe above forces the Amount column to keep being
# stored as text (matching the rest of the sheet) instead of turning into
# a number; drop the resulting "quote prefix" formatting so the cells stay
# on the default (unstyled) style, same as every other cell in the sheet.
$wsExpense.Range("D2:D3").ClearFormats()

# --- Make "Sheet" the active sheet/tab again ---
$wsSheet = $wb.Worksheets.Item("Sheet")
$wsSheet.Activate()
$wsCategory.Range("A1").Select()
$wsExpense.Range("A1").Select()
$wsSheet.Range("A1").Select()

$excel.Calculate()
